# Update the cryptos worksheet with the latest scraped data
# (GitHub Actions scheduled refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28/29 and 32/33 got re-ordered (ranking changed) in addition to
# value refreshes, so we set B/C/D/E for every data row explicitly.
# A leading "'" forces Excel to keep values such as "593.64" / "1.00" as
# plain text (matching the source workbook's inline-string cells) instead
# of re-interpreting them as numbers; re-applying the "Normal" style right
# after keeps the cell formatting untouched.

$updates = @(
    @{ Row=2;  B="Bitcoin";                       C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc";                         D="67.625.50"; E="  -1.58%  " },
    @{ Row=3;  B="Ethereum";                      C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth";                        D="3.785.36";  E="  +0.36%  " },
    @{ Row=4;  B="TetherUSD";                     C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt";                      D="0.999";     E="  -0.19%  " },
    @{ Row=5;  B="BNB";                           C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb";                             D="593.64";    E="  -0.66%  " },
    @{ Row=6;  B="Solana";                        C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol";                              D="166.48";    E="  -1.57%  " },
    @{ Row=7;  B="LidoStakedEther";                C="https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth";                  D="3.788.14";  E="  +0.43%  " },
    @{ Row=8;  B="USDC";                          C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc";                           D="1.00";      E="  +0.06%  " },
    @{ Row=9;  B="XRP";                           C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp";                             D="0.518";     E="  -1.03%  " },
    @{ Row=10; B="Dogecoin";                      C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge";                       D="0.159";     E="  -1.60%  " },
    @{ Row=11; B="Toncoin";                       C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton";                             D="6.36";      E="  -1.97%  " },
    @{ Row=12; B="Cardano";                       C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada";                         D="0.448";     E="  -0.85%  " },
    @{ Row=13; B="ShibaInu";                      C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib";                           D="0.0000255"; E="  -2.54%  " },
    @{ Row=14; B="Avalanche";                     C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax";                          D="35.95";     E="  -1.80%  " },
    @{ Row=15; B="WrappedliquidstakedEther2.0";   C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth";       D="4.418.17";  E="  +0.21%  " },
    @{ Row=16; B="WrappedEther";                  C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth";                   D="3.811.14";  E="  +0.71%  " },
    @{ Row=17; B="WrappedBTC";                    C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc";                       D="67.524.17"; E="  -1.90%  " },
    @{ Row=18; B="Chainlink";                     C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link";                      D="18.18";     E="  +0.56%  " },
    @{ Row=19; B="TRON";                          C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx";                            D="0.112";     E="  +0.15%  " },
    @{ Row=20; B="Polkadot";                      C="https://coinranking.com/coin/25W7FG7om+polkadot-dot";                            D="6.98";      E="  -0.96%  " },
    @{ Row=21; B="Uniswap";                       C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni";                             D="10.24";     E="  -6.19%  " },
    @{ Row=22; B="BitcoinCash";                   C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch";                     D="459.80";    E="  -2.14%  " },
    @{ Row=23; B="Polygon";                       C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic";                        D="0.697";     E="  -1.23%  " },
    @{ Row=24; B="PEPE";                          C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe";                               D="0.0000151"; E="  +3.39%  " },
    @{ Row=25; B="Litecoin";                      C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc";                        D="83.53";     E="  -1.24%  " },
    @{ Row=26; B="Fetch.AI";                      C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet";                         D="2.14";      E="  -3.75%  " },
    @{ Row=27; B="InternetComputer(DFINITY)";     C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";             D="11.86";     E="  -2.37%  " },
    @{ Row=28; B="RenderToken";                   C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";                    D="10.02";     E="  -1.88%  " },
    @{ Row=29; B="Dai";                           C="https://coinranking.com/coin/MoTuySvg7+dai-dai";                                 D="1.00";      E="  +0.10%  " },
    @{ Row=30; B="PancakeSwap";                   C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";                        D="2.78";      E="  -1.27%  " },
    @{ Row=31; B="EthereumClassic";               C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc";                 D="29.85";     E="  -1.08%  " },
    @{ Row=32; B="NEARProtocol";                  C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near";                       D="7.22";      E="  -2.54%  " },
    @{ Row=33; B="ImmutableX";                    C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";                          D="2.19";      E="  -1.11%  " },
    @{ Row=34; B="Aptos";                         C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                               D="9.17";      E="  -1.71%  " },
    @{ Row=35; B="Binance-PegBSC-USD";            C="https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd";              D="0.999";     E="  +0.02%  " },
    @{ Row=36; B="RenzoRestakedETH";              C="https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth";                  D="3.737.09";  E="  +0.20%  " },
    @{ Row=37; B="Hedera";                        C="https://coinranking.com/coin/jad286TjB+hedera-hbar";                             D="0.1000";    E="  -1.75%  " },
    @{ Row=38; B="dogwifhat";                     C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif";                           D="3.32";      E="  -4.97%  " },
    @{ Row=39; B="Kaspa";                         C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas";                               D="0.138";     E="  -0.67%  " },
    @{ Row=40; B="Mantle";                        C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";                              D="0.996";     E="  -0.87%  " },
    @{ Row=41; B="Filecoin";                      C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil";                            D="5.74";      E="  -1.88%  " },
    @{ Row=42; B="FirstDigitalUSD";               C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd";                   D="1.00";      E="  -0.29%  " },
    @{ Row=43; B="USDe";                          C="https://coinranking.com/coin/exbfr2U-0+usde-usde";                               D="1.00";      E="  -0.01%  " },
    @{ Row=44; B="Arweave";                       C="https://coinranking.com/coin/7XWg41D1+arweave-ar";                               D="44.00";     E="  +0.64%  " },
    @{ Row=45; B="TheGraph";                      C="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt";                            D="0.297";     E="  -3.40%  " },
    @{ Row=46; B="OKB";                           C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb";                             D="47.00";     E="  +2.15%  " },
    @{ Row=47; B="Cosmos";                        C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom";                         D="8.37";      E="  -2.96%  " },
    @{ Row=48; B="Monero";                        C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";                          D="147.49";    E="  +0.91%  " },
    @{ Row=49; B="Bittensor";                     C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao";                           D="392.82";    E="  -1.27%  " },
    @{ Row=50; B="Stacks";                        C="https://coinranking.com/coin/mMPrMcB7+stacks-stx";                               D="1.82";      E="  -7.29%  " },
    @{ Row=51; B="Maker";                         C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";                           D="2.756.56";  E="  +2.40%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = "'" + $u.B
    $cellB.Style = "Normal"

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = "'" + $u.C
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = "'" + $u.D
    $cellD.Style = "Normal"

    $cellE = $ws.Cells.Item($r, 5)
    $cellE.Value = "'" + $u.E
    $cellE.Style = "Normal"
}
